# refactoring of xls2ajf function
# Adds a new "select_one MONITORING_VISITS" question plus three follow-up
# "text" questions to the bottom of the "Monitoring visits" group on the
# "survey" sheet (rows 8-11), and moves the group's closing "end group" row
# down to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 8 used to be the lone "end group" row closing the "Monitoring visits"
# group; it now becomes the new select_one question, and the "end group" is
# pushed down to row 12.
$ws.Cells.Item(8, 1).Value = "select_one MONITORING_VISITS"
$ws.Cells.Item(8, 2).Value = "_2_200"
$ws.Cells.Item(8, 3).Value = "Monitoring visits"

$ws.Cells.Item(9, 1).Value = "text"
$ws.Cells.Item(9, 2).Value = "_2_202"
$ws.Cells.Item(9, 3).Value = "Name of monitoring visitor"

$ws.Cells.Item(10, 1).Value = "text"
$ws.Cells.Item(10, 2).Value = "_2_203"
$ws.Cells.Item(10, 3).Value = "Position of monitoring visitor"

$ws.Cells.Item(11, 1).Value = "text"
$ws.Cells.Item(11, 2).Value = "_2_204"
$ws.Cells.Item(11, 3).Value = "Main conclusions of the monitoring visit"

$ws.Cells.Item(12, 1).Value = "end group"

# Column C grew wider to fit the new, longer label text.
$ws.Columns.Item(3).ColumnWidth = 36.26

# Leave the cursor where the author left it after the edit.
$ws.Range("A18").Select()
